$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text, since many values look numeric
# (e.g. "348.28") but must remain stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '52.212.80'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '2.796.20'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '348.28'
$ws.Range("E5").Value = '  +4.55%  '
$ws.Range("D6").Value = '115.83'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '0.549'
$ws.Range("E7").Value = '  +3.51%  '
$ws.Range("D9").Value = '0.589'
$ws.Range("E9").Value = '  +2.22%  '
$ws.Range("D10").Value = '42.39'
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("D11").Value = '0.0860'
$ws.Range("E11").Value = '  +3.94%  '
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '7.85'
$ws.Range("E14").Value = '  +3.15%  '
$ws.Range("D15").Value = '3.238.21'
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("D16").Value = '2.791.27'
$ws.Range("E16").Value = '  +1.80%  '
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '52.201.55'
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D19").Value = '3.17'
$ws.Range("E19").Value = '  +6.31%  '
$ws.Range("E20").Value = '  +6.40%  '
$ws.Range("D21").Value = '13.36'
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("D23").Value = '269.96'
$ws.Range("E23").Value = '  -2.19%  '
$ws.Range("D24").Value = '70.04'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '2.77'
$ws.Range("E25").Value = '  +3.66%  '
$ws.Range("D26").Value = '26.81'
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '10.25'
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("E29").Value = '  +1.06%  '
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").Value = '34.36'
$ws.Range("E31").Value = '  -3.12%  '
$ws.Range("D32").Value = '50.37'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.77'
$ws.Range("E33").Value = '  +2.75%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").Value = '0.0436'
$ws.Range("E34").Value = '  +25.04%  '
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = '18.67'
$ws.Range("E38").Value = '  -3.84%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '4.92'
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("E41").Value = '  +10.02%  '
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("D44").Value = '126.31'
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.056.61'
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").Value = '2.33'
$ws.Range("E48").Value = '  +3.66%  '
$ws.Range("D49").Value = '0.959'
$ws.Range("E49").Value = '  +10.82%  '
$ws.Range("D50").Value = '5.59'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").Value = '8.99'
$ws.Range("E51").Value = '  +0.08%  '
